$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:L1").Font.Bold = $true

$ws.Range("A2").Value = 0.025024199
$ws.Range("B2").Value = 0.2220923
$ws.Range("E2").Value = 0.0534846
$ws.Range("F2").Value = 0.0597214

$ws.Range("A3").Value = 0.0268479
$ws.Range("B3").Value = 0.0130032
$ws.Range("C3").Value = 0.0055264
$ws.Range("D3").Value = 0.0073712
$ws.Range("E3").Value = 0.052746
$ws.Range("F3").Value = 0.0569792

$ws.Range("A4").Value = 0.0079994
$ws.Range("B4").Value = 0.0228975
$ws.Range("C4").Value = 0.0026108
$ws.Range("D4").Value = 0.0025787

$ws.Range("C5").Value = 0.0037162
$ws.Range("D5").Value = 0.0027706

$ws.Range("C15").Select()
